$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 1151; this pushes the existing
# rows 1151-1246 down to 1154-1249 (matching the target dimension A1:R1249).
$ws.Rows("1151:1153").Insert()

# Row 1151 - new "Camote" record for Provincia del Biobío
$ws.Range("A1151").Value = 9
$ws.Range("B1151").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1151").Value = "Metropolitana"
$ws.Range("D1151").Value = 44783
$ws.Range("E1151").Value = 13
$ws.Range("F1151").Value = 100114013
$ws.Range("G1151").Value = "Zanahoria"
$ws.Range("H1151").Value = "Sin especificar"
$ws.Range("I1151").Value = "Camote"
$ws.Range("J1151").Value = 160
$ws.Range("K1151").Value = 8000
$ws.Range("L1151").Value = 8000
$ws.Range("M1151").Value = 8000
$ws.Range("N1151").Value = "$/saco 20 kilos"
$ws.Range("O1151").Value = "Provincia del Biobío"
$ws.Range("P1151").Value = 400
$ws.Range("Q1151").Value = 20
$ws.Range("R1151").Value = "Hortaliza"

# Row 1152 - new "Primera" record for Provincia del Biobío
$ws.Range("A1152").Value = 9
$ws.Range("B1152").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1152").Value = "Metropolitana"
$ws.Range("D1152").Value = 44783
$ws.Range("E1152").Value = 13
$ws.Range("F1152").Value = 100114013
$ws.Range("G1152").Value = "Zanahoria"
$ws.Range("H1152").Value = "Sin especificar"
$ws.Range("I1152").Value = "Primera"
$ws.Range("J1152").Value = 340
$ws.Range("K1152").Value = 10000
$ws.Range("L1152").Value = 10000
$ws.Range("M1152").Value = 10000
$ws.Range("N1152").Value = "$/saco 20 kilos"
$ws.Range("O1152").Value = "Provincia del Biobío"
$ws.Range("P1152").Value = 500
$ws.Range("Q1152").Value = 20
$ws.Range("R1152").Value = "Hortaliza"

# Row 1153 - new "Segunda" record for Provincia del Biobío
$ws.Range("A1153").Value = 9
$ws.Range("B1153").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1153").Value = "Metropolitana"
$ws.Range("D1153").Value = 44783
$ws.Range("E1153").Value = 13
$ws.Range("F1153").Value = 100114013
$ws.Range("G1153").Value = "Zanahoria"
$ws.Range("H1153").Value = "Sin especificar"
$ws.Range("I1153").Value = "Segunda"
$ws.Range("J1153").Value = 250
$ws.Range("K1153").Value = 9000
$ws.Range("L1153").Value = 9000
$ws.Range("M1153").Value = 9000
$ws.Range("N1153").Value = "$/saco 20 kilos"
$ws.Range("O1153").Value = "Provincia del Biobío"
$ws.Range("P1153").Value = 450
$ws.Range("Q1153").Value = 20
$ws.Range("R1153").Value = "Hortaliza"
